$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture the current "yellow highlight" formatting (style index 4:
#     fontId=0, fillId=2 (solid yellow FFFF00), borderId=1) from row 12, A12:C12,
#     and apply it to row 16 (A16:C16) which currently has no fill (style 1). ---
$ws.Range("A12:C12").Copy() | Out-Null
$ws.Range("A16:C16").PasteSpecial(-4122) | Out-Null

# --- Step 2: build the new "white / theme background" fill+border style and
#     apply it to row 12 (A12:C12) - solid fill, fgColor theme Light1 (-> theme="0"),
#     keeping the existing thin border. ---
$ws.Range("A12:C12").Interior.Pattern = 1
$ws.Range("A12:C12").Interior.ThemeColor = 2

# --- Step 3: swap the Y/N values in column C for rows 12 and 16. ---
$ws.Range("C12").Value = "N"
$ws.Range("C16").Value = "Y"

# --- Step 4: update the active selection to match the saved view state. ---
$ws.Range("E5").Select() | Out-Null

$excel.CutCopyMode = 0
